# New crime data collected — weekly CompStat report refresh (cs-en-us-030pct)
#
# Updates:
#   1) Header text: "Volume 32 Number 31" -> "Volume 32 Number 32"
#   2) Header text: report week date range 7/28/2025-8/3/2025 -> 8/4/2025-8/10/2025
#   3) Refreshed weekly/28-day/YTD/2-year crime-count and % change figures
#      for rows 14-33 of the Crime Complaints table.
#
# Several cells in the table flip between a numeric value and the sheet's
# "no data" placeholder text ("0" / "***.*", stored as literal text rather
# than a number) as the underlying counts move to/from zero or undefined
# ratios. Those are handled specially below: the cell is first marked as
# Text so the literal digits aren't re-parsed as a number, then its
# formatting is replaced wholesale (via copy/paste of an untouched donor
# cell's format) so the resulting style index exactly matches a plain
# numeric-column cell rather than a newly-minted "text number format" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Volume / Number header -------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  32"

# --- 2) Report-covering-the-week header ----------------------------------
$ws.Range("C9").Value = "Report Covering the Week  8/4/2025  Through  8/10/2025"

# --- 3) Plain numeric value refreshes (style/type unchanged) ------------
$ws.Range("N14").Value = -97.297297297297
$ws.Range("F15").Value = 1
$ws.Range("L15").Value = 100
$ws.Range("N15").Value = -76.470588235294
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 78
$ws.Range("K16").Value = -24.358974358974
$ws.Range("L16").Value = -22.368421052631
$ws.Range("M16").Value = -54.615384615384
$ws.Range("N16").Value = -86.150234741784
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -21.052631578947
$ws.Range("I17").Value = 101
$ws.Range("J17").Value = 143
$ws.Range("K17").Value = -29.370629370629
$ws.Range("L17").Value = -28.873239436619
$ws.Range("M17").Value = -12.173913043478
$ws.Range("N17").Value = -77.704194260485
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = -25.423728813559
$ws.Range("L18").Value = 7.317073170731
$ws.Range("M18").Value = -25.423728813559
$ws.Range("N18").Value = -90.393013100436
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 20
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 221
$ws.Range("J19").Value = 215
$ws.Range("K19").Value = 2.790697674418
$ws.Range("L19").Value = 8.333333333333
$ws.Range("M19").Value = 148.314606741573
$ws.Range("N19").Value = 12.182741116751
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 14.634146341463
$ws.Range("L20").Value = -32.857142857142
$ws.Range("M20").Value = 51.612903225806
$ws.Range("N20").Value = -74.866310160427
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -27.777777777777
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -13.846153846153
$ws.Range("I21").Value = 481
$ws.Range("J21").Value = 545
$ws.Range("K21").Value = -11.743119266055
$ws.Range("L21").Value = -10.925925925925
$ws.Range("M21").Value = 7.366071428571
$ws.Range("N21").Value = -73.158482142857
$ws.Range("L22").Value = -26.315789473684
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = -20
$ws.Range("I24").Value = 476
$ws.Range("J24").Value = 466
$ws.Range("K24").Value = 2.145922746781
$ws.Range("L24").Value = -5.367793240556
$ws.Range("M24").Value = 118.348623853211
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = -64
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 127
$ws.Range("K25").Value = -27.559055118110
$ws.Range("L25").Value = -31.343283582089
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -21.428571428571
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 29.629629629629
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 210
$ws.Range("K26").Value = 4.761904761904
$ws.Range("L26").Value = 16.402116402116
$ws.Range("M26").Value = -27.631578947368
$ws.Range("F27").Value = 1
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 3
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 5

# --- 4) Cells that flip from a number to literal placeholder text -------
# (format copied from row 14's untouched style-13 donors so the
# resulting style/type exactly matches the rest of the "no data" cells)

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G28").PasteSpecial(-4122)

$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H28").PasteSpecial(-4122)

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# --- 5) Cell that flips from placeholder text back to a number ----------
$ws.Range("L33").Value = 0
$ws.Range("K14").Copy()
$ws.Range("L33").PasteSpecial(-4122)
